$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'275.04"
$ws.Range("D3").Formula = "'22.87"
$ws.Range("D4").Formula = "'6.327"
$ws.Range("D5").Formula = "'0.06236"
$ws.Range("D6").Formula = "'3.649"
$ws.Range("D7").Formula = "'6.639"
$ws.Range("D8").Formula = "'1.399"
$ws.Range("D9").Formula = "'0.8336"
$ws.Range("D10").Formula = "'0.01381"
$ws.Range("D11").Formula = "'0.1605"
$ws.Range("D12").Formula = "'0.08387"
$ws.Range("D13").Formula = "'0.03541"
$ws.Range("D14").Formula = "'0.03181"
$ws.Range("D15").Formula = "'4.086"
$ws.Range("D16").Formula = "'0.09292"
$ws.Range("D17").Formula = "'0.001673"
$ws.Range("D18").Formula = "'0.04747"
$ws.Range("D19").Formula = "'0.006396"
$ws.Range("D20").Formula = "'0.005714"
$ws.Range("E20").Value = "19HotbitTokenHTBWorstin24h"
$ws.Range("D23").Formula = "'3.722"
$ws.Range("D24").Formula = "'2.326"
$ws.Range("D25").Formula = "'0.3326"
$ws.Range("D28").Formula = "'0.0002707"
$ws.Range("D40").Formula = "'0.04730"
$ws.Range("D41").Formula = "'0.007096"
$ws.Range("D42").Formula = "'0.1170"
$ws.Range("D43").Formula = "'0.003656"
$ws.Range("D44").Formula = "'0.01217"
$ws.Range("D45").Formula = "'0.00006032"
$ws.Range("D46").Formula = "'0.0009911"
$ws.Range("D47").Formula = "'0.00000000751"
$ws.Range("D48").Formula = "'0.7830"
$ws.Range("D49").Formula = "'0.002421"
$ws.Range("D50").Formula = "'0.00002402"
$ws.Range("E50").Value = "49CryptobidCoinCBC"
$ws.Range("D51").Formula = "'0.01241"
